$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" status text to "In Translation" on each sheet.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Shrink the "Status" columns that previously held the long "Ready for handoff" text.
# (ColumnWidth snaps to the nearest 1/6-character pixel grid, so 12.5 is the input
# that lands closest to the target stored width of ~13.41.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
